$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-09-23 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-24 Sunday", 2) | Out-Null
$d.Content.Find.Execute("63+12=75", $true, $false, $false, $false, $false, $true, 1, $false, "83-82=1", 2) | Out-Null
$d.Content.Find.Execute("86-24=62", $true, $false, $false, $false, $false, $true, 1, $false, "71+18=89", 2) | Out-Null
$d.Content.Find.Execute("47-27=20", $true, $false, $false, $false, $false, $true, 1, $false, "81-11=70", 2) | Out-Null
$d.Content.Find.Execute("82-49=33", $true, $false, $false, $false, $false, $true, 1, $false, "35-5=30", 2) | Out-Null
$d.Content.Find.Execute("3+48=51", $true, $false, $false, $false, $false, $true, 1, $false, "27+34=61", 2) | Out-Null
$d.Content.Find.Execute("73-47=26", $true, $false, $false, $false, $false, $true, 1, $false, "50-27=23", 2) | Out-Null
$d.Content.Find.Execute("2+7=9", $true, $false, $false, $false, $false, $true, 1, $false, "76+4=80", 2) | Out-Null
$d.Content.Find.Execute("54+10=64", $true, $false, $false, $false, $false, $true, 1, $false, "36+41=77", 2) | Out-Null
$d.Content.Find.Execute("85-27=58", $true, $false, $false, $false, $false, $true, 1, $false, "2+25=27", 2) | Out-Null
$d.Content.Find.Execute("80-76=4", $true, $false, $false, $false, $false, $true, 1, $false, "57+34=91", 2) | Out-Null
$d.Content.Find.Execute("17-13=4", $true, $false, $false, $false, $false, $true, 1, $false, "33+62=95", 2) | Out-Null
$d.Content.Find.Execute("49+35=84", $true, $false, $false, $false, $false, $true, 1, $false, "81-57=24", 2) | Out-Null
$d.Content.Find.Execute("68-63=5", $true, $false, $false, $false, $false, $true, 1, $false, "14+49=63", 2) | Out-Null
$d.Content.Find.Execute("62-9=53", $true, $false, $false, $false, $false, $true, 1, $false, "55-0=55", 2) | Out-Null
$d.Content.Find.Execute("61-55=6", $true, $false, $false, $false, $false, $true, 1, $false, "28+67=95", 2) | Out-Null
$d.Content.Find.Execute("29+0=29", $true, $false, $false, $false, $false, $true, 1, $false, "10+26=36", 2) | Out-Null
$d.Content.Find.Execute("68-35=33", $true, $false, $false, $false, $false, $true, 1, $false, "67-49=18", 2) | Out-Null
$d.Content.Find.Execute("15-6=9", $true, $false, $false, $false, $false, $true, 1, $false, "32+58=90", 2) | Out-Null
$d.Content.Find.Execute("75-57=18", $true, $false, $false, $false, $false, $true, 1, $false, "70-27=43", 2) | Out-Null
$d.Content.Find.Execute("66-7=59", $true, $false, $false, $false, $false, $true, 1, $false, "51-48=3", 2) | Out-Null
$d.Content.Find.Execute("13+22=35", $true, $false, $false, $false, $false, $true, 1, $false, "16+17=33", 2) | Out-Null
$d.Content.Find.Execute("28+71=99", $true, $false, $false, $false, $false, $true, 1, $false, "29+58=87", 2) | Out-Null
$d.Content.Find.Execute("11+54=65", $true, $false, $false, $false, $false, $true, 1, $false, "95-70=25", 2) | Out-Null
$d.Content.Find.Execute("95+0=95", $true, $false, $false, $false, $false, $true, 1, $false, "29+20=49", 2) | Out-Null
$d.Content.Find.Execute("57-12=45", $true, $false, $false, $false, $false, $true, 1, $false, "98-40=58", 2) | Out-Null
$d.Content.Find.Execute("24-15=9", $true, $false, $false, $false, $false, $true, 1, $false, "49+12=61", 2) | Out-Null
$d.Content.Find.Execute("62+26=88", $true, $false, $false, $false, $false, $true, 1, $false, "15+23=38", 2) | Out-Null
$d.Content.Find.Execute("35+22=57", $true, $false, $false, $false, $false, $true, 1, $false, "2+63=65", 2) | Out-Null
$d.Content.Find.Execute("79-76=3", $true, $false, $false, $false, $false, $true, 1, $false, "39-33=6", 2) | Out-Null
$d.Content.Find.Execute("45+8=53", $true, $false, $false, $false, $false, $true, 1, $false, "17-10=7", 2) | Out-Null
$d.Content.Find.Execute("94-6=88", $true, $false, $false, $false, $false, $true, 1, $false, "49+11=60", 2) | Out-Null
$d.Content.Find.Execute("89-74=15", $true, $false, $false, $false, $false, $true, 1, $false, "90-3=87", 2) | Out-Null
$d.Content.Find.Execute("13+33=46", $true, $false, $false, $false, $false, $true, 1, $false, "16-5=11", 2) | Out-Null
$d.Content.Find.Execute("9+8=17", $true, $false, $false, $false, $false, $true, 1, $false, "51-48=3", 2) | Out-Null
$d.Content.Find.Execute("25+73=98", $true, $false, $false, $false, $false, $true, 1, $false, "93-36=57", 2) | Out-Null
$d.Content.Find.Execute("92-57=35", $true, $false, $false, $false, $false, $true, 1, $false, "97-60=37", 2) | Out-Null
$d.Content.Find.Execute("37+25=62", $true, $false, $false, $false, $false, $true, 1, $false, "83-45=38", 2) | Out-Null
$d.Content.Find.Execute("10+64=74", $true, $false, $false, $false, $false, $true, 1, $false, "49+41=90", 2) | Out-Null
$d.Content.Find.Execute("70-64=6", $true, $false, $false, $false, $false, $true, 1, $false, "93-7=86", 2) | Out-Null
$d.Content.Find.Execute("84-23=61", $true, $false, $false, $false, $false, $true, 1, $false, "93-79=14", 2) | Out-Null
$d.Content.Find.Execute("67-27=40", $true, $false, $false, $false, $false, $true, 1, $false, "82-62=20", 2) | Out-Null
$d.Content.Find.Execute("31+34=65", $true, $false, $false, $false, $false, $true, 1, $false, "69-40=29", 2) | Out-Null
$d.Content.Find.Execute("97-4=93", $true, $false, $false, $false, $false, $true, 1, $false, "14+46=60", 2) | Out-Null
$d.Content.Find.Execute("2-0=2", $true, $false, $false, $false, $false, $true, 1, $false, "57-52=5", 2) | Out-Null
$d.Content.Find.Execute("98-70=28", $true, $false, $false, $false, $false, $true, 1, $false, "2+18=20", 2) | Out-Null
$d.Content.Find.Execute("52+18=70", $true, $false, $false, $false, $false, $true, 1, $false, "25+43=68", 2) | Out-Null
$d.Content.Find.Execute("60+20=80", $true, $false, $false, $false, $false, $true, 1, $false, "73+2=75", 2) | Out-Null
$d.Content.Find.Execute("81-18=63", $true, $false, $false, $false, $false, $true, 1, $false, "79-70=9", 2) | Out-Null
$d.Content.Find.Execute("1+35=36", $true, $false, $false, $false, $false, $true, 1, $false, "24+4=28", 2) | Out-Null
$d.Content.Find.Execute("99-81=18", $true, $false, $false, $false, $false, $true, 1, $false, "53-9=44", 2) | Out-Null
$d.Content.Find.Execute("67-57=10", $true, $false, $false, $false, $false, $true, 1, $false, "96-52=44", 2) | Out-Null
$d.Content.Find.Execute("41+36=77", $true, $false, $false, $false, $false, $true, 1, $false, "44-36=8", 2) | Out-Null
$d.Content.Find.Execute("51+39=90", $true, $false, $false, $false, $false, $true, 1, $false, "87-38=49", 2) | Out-Null
$d.Content.Find.Execute("67-23=44", $true, $false, $false, $false, $false, $true, 1, $false, "27+21=48", 2) | Out-Null
$d.Content.Find.Execute("9+78=87", $true, $false, $false, $false, $false, $true, 1, $false, "39+1=40", 2) | Out-Null
$d.Content.Find.Execute("76+20=96", $true, $false, $false, $false, $false, $true, 1, $false, "11-7=4", 2) | Out-Null
$d.Content.Find.Execute("24+5=29", $true, $false, $false, $false, $false, $true, 1, $false, "86-81=5", 2) | Out-Null
$d.Content.Find.Execute("43+21=64", $true, $false, $false, $false, $false, $true, 1, $false, "72+23=95", 2) | Out-Null
$d.Content.Find.Execute("47+12=59", $true, $false, $false, $false, $false, $true, 1, $false, "34+59=93", 2) | Out-Null
$d.Content.Find.Execute("10+38=48", $true, $false, $false, $false, $false, $true, 1, $false, "8+41=49", 2) | Out-Null
$d.Content.Find.Execute("33+2=35", $true, $false, $false, $false, $false, $true, 1, $false, "34-8=26", 2) | Out-Null
$d.Content.Find.Execute("91-3=88", $true, $false, $false, $false, $false, $true, 1, $false, "87-62=25", 2) | Out-Null
$d.Content.Find.Execute("26+41=67", $true, $false, $false, $false, $false, $true, 1, $false, "39+46=85", 2) | Out-Null
$d.Content.Find.Execute("24+44=68", $true, $false, $false, $false, $false, $true, 1, $false, "63+9=72", 2) | Out-Null
$d.Content.Find.Execute("57+9=66", $true, $false, $false, $false, $false, $true, 1, $false, "59-54=5", 2) | Out-Null
$d.Content.Find.Execute("20+20=40", $true, $false, $false, $false, $false, $true, 1, $false, "70-59=11", 2) | Out-Null
$d.Content.Find.Execute("32+54=86", $true, $false, $false, $false, $false, $true, 1, $false, "78+16=94", 2) | Out-Null
$d.Content.Find.Execute("80-79=1", $true, $false, $false, $false, $false, $true, 1, $false, "40-32=8", 2) | Out-Null
$d.Content.Find.Execute("40-1=39", $true, $false, $false, $false, $false, $true, 1, $false, "84-51=33", 2) | Out-Null
$d.Content.Find.Execute("13+10=23", $true, $false, $false, $false, $false, $true, 1, $false, "10+19=29", 2) | Out-Null
$d.Content.Find.Execute("44+10=54", $true, $false, $false, $false, $false, $true, 1, $false, "93-32=61", 2) | Out-Null
$d.Content.Find.Execute("28+31=59", $true, $false, $false, $false, $false, $true, 1, $false, "63-2=61", 2) | Out-Null
$d.Content.Find.Execute("83+1=84", $true, $false, $false, $false, $false, $true, 1, $false, "58+18=76", 2) | Out-Null
$d.Content.Find.Execute("16+36=52", $true, $false, $false, $false, $false, $true, 1, $false, "78-14=64", 2) | Out-Null
$d.Content.Find.Execute("17+9=26", $true, $false, $false, $false, $false, $true, 1, $false, "2+62=64", 2) | Out-Null
$d.Content.Find.Execute("23+0=23", $true, $false, $false, $false, $false, $true, 1, $false, "39-8=31", 2) | Out-Null
$d.Content.Find.Execute("80-4=76", $true, $false, $false, $false, $false, $true, 1, $false, "21+37=58", 2) | Out-Null
$d.Content.Find.Execute("12-4=8", $true, $false, $false, $false, $false, $true, 1, $false, "96-92=4", 2) | Out-Null
$d.Content.Find.Execute("63-11=52", $true, $false, $false, $false, $false, $true, 1, $false, "46+43=89", 2) | Out-Null
$d.Content.Find.Execute("19+21=40", $true, $false, $false, $false, $false, $true, 1, $false, "36+2=38", 2) | Out-Null
$d.Content.Find.Execute("93-34=59", $true, $false, $false, $false, $false, $true, 1, $false, "57-36=21", 2) | Out-Null
$d.Content.Find.Execute("48+30=78", $true, $false, $false, $false, $false, $true, 1, $false, "63-10=53", 2) | Out-Null
$d.Content.Find.Execute("12+86=98", $true, $false, $false, $false, $false, $true, 1, $false, "93-23=70", 2) | Out-Null
$d.Content.Find.Execute("59-34=25", $true, $false, $false, $false, $false, $true, 1, $false, "67+2=69", 2) | Out-Null
$d.Content.Find.Execute("4+57=61", $true, $false, $false, $false, $false, $true, 1, $false, "66+9=75", 2) | Out-Null
$d.Content.Find.Execute("66-52=14", $true, $false, $false, $false, $false, $true, 1, $false, "43-33=10", 2) | Out-Null
$d.Content.Find.Execute("25+14=39", $true, $false, $false, $false, $false, $true, 1, $false, "64-27=37", 2) | Out-Null
$d.Content.Find.Execute("7+63=70", $true, $false, $false, $false, $false, $true, 1, $false, "10-8=2", 2) | Out-Null
$d.Content.Find.Execute("74-56=18", $true, $false, $false, $false, $false, $true, 1, $false, "45-45=0", 2) | Out-Null
$d.Content.Find.Execute("99-21=78", $true, $false, $false, $false, $false, $true, 1, $false, "92-76=16", 2) | Out-Null
$d.Content.Find.Execute("27-19=8", $true, $false, $false, $false, $false, $true, 1, $false, "5+53=58", 2) | Out-Null
$d.Content.Find.Execute("45+46=91", $true, $false, $false, $false, $false, $true, 1, $false, "91-89=2", 2) | Out-Null
$d.Content.Find.Execute("85-9=76", $true, $false, $false, $false, $false, $true, 1, $false, "79-54=25", 2) | Out-Null
$d.Content.Find.Execute("60-1=59", $true, $false, $false, $false, $false, $true, 1, $false, "44-0=44", 2) | Out-Null
$d.Content.Find.Execute("34+6=40", $true, $false, $false, $false, $false, $true, 1, $false, "4+66=70", 2) | Out-Null
$d.Content.Find.Execute("91-42=49", $true, $false, $false, $false, $false, $true, 1, $false, "98-34=64", 2) | Out-Null
$d.Content.Find.Execute("76-43=33", $true, $false, $false, $false, $false, $true, 1, $false, "15+7=22", 2) | Out-Null
$d.Content.Find.Execute("17+69=86", $true, $false, $false, $false, $false, $true, 1, $false, "25+28=53", 2) | Out-Null
$d.Content.Find.Execute("54+0=54", $true, $false, $false, $false, $false, $true, 1, $false, "17-2=15", 2) | Out-Null
$d.Content.Find.Execute("84-4=80", $true, $false, $false, $false, $false, $true, 1, $false, "22+41=63", 2) | Out-Null

Write-Host "Replacements applied: 101"
